# IGCC Netting Flows Historical - roll the data window forward by 2 days
# (retraining the model for Astro). Column A holds the 15-minute timestamps,
# column D holds the quarter-of-day index (1..96, unchanged), and column E
# holds a "Lookup" string built from the date portion of column A plus the
# quarter index in column D (e.g. "21.01.20261"). This script shifts every
# timestamp in column A forward by 2 calendar days and regenerates the
# corresponding Lookup strings in column E to match, leaving columns B, C
# and D untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1
$headerRow = 1
$dataFirstRow = $headerRow + 1

if ($lastRow -lt $dataFirstRow) {
    Write-Host "No data rows found."
} else {
    $n = $lastRow - $dataFirstRow + 1

    $rngA = $ws.Range($ws.Cells.Item($dataFirstRow, 1), $ws.Cells.Item($lastRow, 1))
    $rngD = $ws.Range($ws.Cells.Item($dataFirstRow, 4), $ws.Cells.Item($lastRow, 4))
    $rngE = $ws.Range($ws.Cells.Item($dataFirstRow, 5), $ws.Cells.Item($lastRow, 5))

    $valsA = $rngA.Value2
    $valsD = $rngD.Value2
    $valsE = $rngE.Value2

    $epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

    for ($i = 1; $i -le $n; $i++) {
        $newSerial = $valsA[$i, 1] + 2
        $valsA[$i, 1] = $newSerial

        $dt = $epoch.AddDays($newSerial)
        $datePart = $dt.ToString("dd.MM.yyyy")
        $quarter = [int]$valsD[$i, 1]
        $valsE[$i, 1] = $datePart + $quarter.ToString()
    }

    $rngA.Value2 = $valsA
    $rngE.Value2 = $valsE

    Write-Host "Shifted $n rows forward by 2 days."
}
